$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.248.58"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.489.21"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.43%  "
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.81%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "2.881.76"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "2.489.74"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "47.177.24"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +15.25%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "245.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.85%  "
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("E30").Value = "  +3.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  +2.59%  "
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.26%  "
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("D45").Value = "1.995.83"
$ws.Range("E45").Value = "  +2.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.97%  "
$ws.Range("E47").Value = "  -4.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.36%  "
